$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new record appended below the existing header (row 1) and data (row 2)
$ws.Range("A3").Value = 112453142
$ws.Range("B3").Value = 56575
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 103021
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"

# Empty text cells (present but blank in the source data) - use a quote-prefix
# entry to force a blank text cell, then reset the style back to Normal so no
# stray formatting is left behind.
$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = "Normal"
$ws.Range("K3").Value = "'"
$ws.Range("K3").Style = "Normal"

$ws.Range("P3").Value = "Valsjöbäcken, Jmt"
$ws.Range("Q3").Value = 443414
$ws.Range("R3").Value = 7021456
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Åre"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Mattmar"

# Date/time columns are stored as plain text in this sheet, not real dates, so
# force text formatting before assigning to avoid Excel's automatic date
# parsing, then reset the style back to Normal.
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-10-01"
$ws.Range("Y3").Style = "Normal"

$ws.Range("Z3").Value = "18:02"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-10-01"
$ws.Range("AA3").Style = "Normal"

$ws.Range("AB3").Value = "18:02"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

$ws.Range("AT3").Value = "'"
$ws.Range("AT3").Style = "Normal"

$ws.Range("AW3").Value = "Henrik Berggren"
$ws.Range("AX3").Value = "Henrik Berggren"

$ws.Range("AY3").Value = "'"
$ws.Range("AY3").Style = "Normal"
